# "Generate Report for handback" — mark the handoff packages as handed
# back (in sync with en-US) and record the handback target/file/datetime
# for each localized language sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$mdFile  = "3bac9dce-e042-4ff9-be51-261a1d015ce5.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/4ef6798426261de1dfeb6f48389f954ff47250e0/e2e/3bac9dce-e042-4ff9-be51-261a1d015ce5.md"

$languages = @(
    @{ Sheet = "zh-cn"; XlfFile = "3bac9dce-e042-4ff9-be51-261a1d015ce5.0343cc1b35dc03dcaf1978a748210400f6a412f7.zh-cn.xlf"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0e3eea29b30dd220d9bb62992b752ce100adbbb8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/3bac9dce-e042-4ff9-be51-261a1d015ce5.0343cc1b35dc03dcaf1978a748210400f6a412f7.zh-cn.xlf"; HandbackTime = "2016-02-15 08:46:02" },
    @{ Sheet = "de-de"; XlfFile = "3bac9dce-e042-4ff9-be51-261a1d015ce5.0343cc1b35dc03dcaf1978a748210400f6a412f7.de-de.xlf"; XlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2e3afa56a0f6440659cab813c4f788b1052e77c0/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/3bac9dce-e042-4ff9-be51-261a1d015ce5.0343cc1b35dc03dcaf1978a748210400f6a412f7.de-de.xlf"; HandbackTime = "2016-02-15 08:46:29" }
)

$handedBack = "Handed back: in sync with en-US"

# Overview sheet: column B is zh-cn status, column C is de-de status.
# Both file rows (2 and 3) move from "Ready for handoff" to handed-back.
$overview = $wb.Worksheets.Item("Overview")
$overview.Cells.Item(2, 2).Value = $handedBack
$overview.Cells.Item(2, 3).Value = $handedBack
$overview.Cells.Item(3, 2).Value = $handedBack
$overview.Cells.Item(3, 3).Value = $handedBack

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Rows 2 and 3 both hold handoff packages awaiting handback (row 4 is
    # the ignored .localization-config entry and does not change).
    foreach ($row in 2, 3) {
        # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
        $ws.Cells.Item($row, 2).Value = $handedBack

        # Latest Target File (E) / Latest Handback File (F): record the
        # handed-back source + xlf, each as a hyperlink like the existing
        # Source File Name / Latest Handoff File columns.
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $mdUrl, "", "", $mdFile) | Out-Null
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $lang.XlfUrl, "", "", $lang.XlfFile) | Out-Null

        # Latest Handback DateTime (G)
        $ws.Cells.Item($row, 7).Value = $lang.HandbackTime
    }
}
